$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
# Old layout: A=Nombre B=Correo electrónico C=ID D=Localización E=Tipo
# New layout: A=Nombre B=Correo electrónico C=ID D=Latitud E=Longitud F=Tipo
$ws.Range("F1").Value = "Tipo"
$ws.Range("D1").Value = "Latitud"
$ws.Range("E1").Value = "Longitud"

# --- Row 2 (existing data row) ---
# Old layout: A2=Sensor 1234 B2=sensor1234@example.com(hyperlink) C2=681356515 E2=3
# New layout: A2=Sensor 1234 B2=sensor1234@example.com(hyperlink) C2=681356515 D2=lat E2=long F2=3
$ws.Range("F2").Value = 3
$ws.Range("D2").Value = "lat"
$ws.Range("E2").Value = "long"

# --- Row 3 (new data row) ---
$ws.Range("A3").Value = "Sensor 1234"
$ws.Range("B3").Value = "sensor1234@example.com"
$ws.Range("C3").Value = 681356515
$ws.Range("D3").Value = 15.65
$ws.Range("E3").Value = "long"
$ws.Range("F3").Value = 3

# --- Row 4 (new data row) ---
$ws.Range("A4").Value = "Sensor 1234"
$ws.Range("B4").Value = "sensor1234@example.com"
$ws.Range("C4").Value = 681356515
$ws.Range("D4").Value = "lat"
$ws.Range("E4").Value = 15.65
$ws.Range("F4").Value = 3

# --- Hyperlinks for the new rows' email cells ---
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:sensor1234@example.com")
$ws.Range("B3").Style = "Hipervínculo"

$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:sensor1234@example.com")
$ws.Range("B4").Style = "Hipervínculo"

# --- Selection matches the final active cell in the source workbook ---
$ws.Range("F4").Select()
